# Update cryptocurrency Price (D) and Volume(1h) (E) columns per the
# automated GitHub Actions data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.249.91'
$ws.Range('E2').Value = '  +6.27%  '
$ws.Range('D3').Value = '3.120.02'
$ws.Range('E3').Value = '  +4.21%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'584.90"
$ws.Range('E5').Value = '  +3.68%  '
$ws.Range('D6').Value = "'144.90"
$ws.Range('E6').Value = '  +4.33%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.113.00'
$ws.Range('E8').Value = '  +4.35%  '
$ws.Range('E9').Value = '  +1.71%  '
$ws.Range('E10').Value = '  +12.50%  '
$ws.Range('D11').Value = "'5.80"
$ws.Range('E11').Value = '  +10.03%  '
$ws.Range('E12').Value = '  +3.07%  '
$ws.Range('E13').Value = '  +7.72%  '
$ws.Range('D14').Value = "'35.60"
$ws.Range('E14').Value = '  +5.20%  '
$ws.Range('D16').Value = '3.635.44'
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('D18').Value = '63.171.02'
$ws.Range('E18').Value = '  +6.12%  '
$ws.Range('D19').Value = '3.118.08'
$ws.Range('E19').Value = '  +4.24%  '
$ws.Range('D20').Value = "'467.36"
$ws.Range('E20').Value = '  +7.26%  '
$ws.Range('D21').Value = "'14.08"
$ws.Range('E21').Value = '  +3.72%  '
$ws.Range('D22').Value = "'0.725"
$ws.Range('E22').Value = '  +1.02%  '
$ws.Range('D23').Value = "'7.56"
$ws.Range('E23').Value = '  +6.95%  '
$ws.Range('D24').Value = "'13.29"
$ws.Range('E24').Value = '  -1.54%  '
$ws.Range('E25').Value = '  +2.29%  '
$ws.Range('D27').Value = "'2.23"
$ws.Range('E28').Value = '  +7.44%  '
$ws.Range('E29').Value = '  +5.43%  '
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('E31').Value = '  +10.08%  '
$ws.Range('D32').Value = "'26.94"
$ws.Range('E32').Value = '  +4.48%  '
$ws.Range('E33').Value = '  +3.75%  '
$ws.Range('D34').Value = '0.0₃0863'
$ws.Range('E34').Value = '  +10.93%  '
$ws.Range('D35').Value = "'2.42"
$ws.Range('E35').Value = '  +15.89%  '
$ws.Range('E36').Value = '  +5.44%  '
$ws.Range('D37').Value = "'3.30"
$ws.Range('E37').Value = '  +19.13%  '
$ws.Range('E38').Value = '  +2.32%  '
$ws.Range('D39').Value = "'51.00"
$ws.Range('E39').Value = '  +4.22%  '
$ws.Range('D40').Value = "'432.18"
$ws.Range('E40').Value = '  +7.29%  '
$ws.Range('D41').Value = "'8.72"
$ws.Range('E41').Value = '  +1.59%  '
$ws.Range('D42').Value = '2.927.83'
$ws.Range('E42').Value = '  +6.04%  '
$ws.Range('E43').Value = '  +4.56%  '
$ws.Range('D44').Value = "'0.278"
$ws.Range('E44').Value = '  +11.02%  '
$ws.Range('E45').Value = '  +5.75%  '
$ws.Range('E46').Value = '  +7.65%  '
$ws.Range('D47').Value = "'35.20"
$ws.Range('E47').Value = '  +1.85%  '
$ws.Range('D49').Value = "'123.39"
$ws.Range('E49').Value = '  +0.26%  '
$ws.Range('E50').Value = '  +0.84%  '
$ws.Range('D51').Value = "'24.55"
$ws.Range('E51').Value = '  +4.18%  '
